$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.185015559196472
$ws.Range("B1").Value = 2.441283941268921
$ws.Range("D1").Value = 2.308252334594727
$ws.Range("E1").Value = 1.186413645744324
